# Add season record columns (Wins / Losses / Ties) to the right of the
# existing "Unnamed: 28" column (AC), for each player row.
#
# Per the commit message, the author added a helper that fetches each
# team's season record (wins/losses/ties) and merges it onto every row of
# the per-team roster table. In this workbook every player row gets the
# same team record: 80 wins, 82 losses, 0 ties.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): AD1="Wins", AE1="Losses", AF1="Ties" ---
# Copy the formatting (bold, bordered, centered) from the neighboring
# header cell AC1 so the new headers look consistent with the rest of
# the header row, then overwrite with the new labels.
$ws.Range("AC1").Copy($ws.Range("AD1"))
$ws.Range("AC1").Copy($ws.Range("AE1"))
$ws.Range("AC1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-53): AD=Wins(80), AE=Losses(82), AF=Ties(0) ---
$firstRow = 2
$lastRow = 53

$winsRange = $ws.Range("AD$firstRow`:AD$lastRow")
$lossesRange = $ws.Range("AE$firstRow`:AE$lastRow")
$tiesRange = $ws.Range("AF$firstRow`:AF$lastRow")

$winsRange.Value = 80
$lossesRange.Value = 82
$tiesRange.Value = 0
